$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove obsolete rows 6 and 7 (MuSCs target-cluster rows), shifting the
# dimension down to A1:T5 and dropping the now-unused "MuSCs" shared string.
$ws.Rows("6:7").Delete() | Out-Null

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Icam5"
$ws.Range("C2").Value = "Itgb2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1421496666666667
$ws.Range("H2").Value = 0.426449
$ws.Range("I2").Value = 0.1211014306728536
$ws.Range("J2").Value = 0.1211014306728536
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.05619066666666667
$ws.Range("N2").Value = 0.168572
$ws.Range("O2").Value = 0.3931387525216601
$ws.Range("P2").Value = 0.39313875252166
$ws.Range("Q2").Value = 0.007987484536444446
$ws.Range("R2").Value = 0.071887360828
$ws.Range("S2").Value = 0.04760966538331396
$ws.Range("T2").Value = 0.04760966538331395

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Icam5"
$ws.Range("C3").Value = "Itgb2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1421496666666667
$ws.Range("H3").Value = 0.426449
$ws.Range("I3").Value = 0.1211014306728536
$ws.Range("J3").Value = 0.1211014306728536
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.08673766666666667
$ws.Range("N3").Value = 0.260213
$ws.Range("O3").Value = 0.60686124747834
$ws.Range("P3").Value = 0.60686124747834
$ws.Range("Q3").Value = 0.01232973040411111
$ws.Range("R3").Value = 0.110967573637
$ws.Range("S3").Value = 0.07349176528953963
$ws.Range("T3").Value = 0.07349176528953963

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Icam5"
$ws.Range("C4").Value = "Itgb2"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.031657
$ws.Range("H4").Value = 3.094971
$ws.Range("I4").Value = 0.8788985693271465
$ws.Range("J4").Value = 0.8788985693271465
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.05619066666666667
$ws.Range("N4").Value = 0.168572
$ws.Range("O4").Value = 0.3931387525216601
$ws.Range("P4").Value = 0.39313875252166
$ws.Range("Q4").Value = 0.05796949460133333
$ws.Range("R4").Value = 0.521725451412
$ws.Range("S4").Value = 0.3455290871383462
$ws.Range("T4").Value = 0.3455290871383461

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Icam5"
$ws.Range("C5").Value = "Itgb2"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.031657
$ws.Range("H5").Value = 3.094971
$ws.Range("I5").Value = 0.8788985693271465
$ws.Range("J5").Value = 0.8788985693271465
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.08673766666666667
$ws.Range("N5").Value = 0.260213
$ws.Range("O5").Value = 0.60686124747834
$ws.Range("P5").Value = 0.60686124747834
$ws.Range("Q5").Value = 0.08948352098033334
$ws.Range("R5").Value = 0.8053516888230001
$ws.Range("S5").Value = 0.5333694821888004
$ws.Range("T5").Value = 0.5333694821888004
